$wb = $excel.ActiveWorkbook

# --- Update selections on the existing sheets (matches author's click-through) ---
$wsNum = $wb.Worksheets.Item("Numero spettacoli")
$wsNum.Activate()
$wsNum.Range("B2:F13").Select()

$wsIng = $wb.Worksheets.Item("Ingressi")
$wsIng.Activate()
$wsIng.Range("B2:F13").Select()

$wsSpe = $wb.Worksheets.Item("Spesa del pubblico")
$wsSpe.Activate()
$wsSpe.Range("B2:F13").Select()

# --- Add the new "Sheet1" worksheet at the end of the workbook ---
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "Sheet1"
$wsNew.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNew = $wb.Worksheets.Item("Sheet1")

# --- Header row ---
$wsNew.Range("A1").Value = "Date"
$wsNew.Range("B1").Value = "Numero spettacoli"
$wsNew.Range("C1").Value = "Ingressi"
$wsNew.Range("D1").Value = "Spesa del pubblico"

# --- Monthly data rows (Jan 2018 - Dec 2022) ---
$data = @(
  @(43101, 407032, 23337474, 482597644.24999994),
  @(43132, 375336, 22230691, 301164981.49999994),
  @(43160, 407807, 21244449, 328709780.68999994),
  @(43191, 388672, 21888769, 461044466.70000011),
  @(43221, 356118, 16985095, 355097932.00999856),
  @(43252, 315634, 14359139, 471409667.89999866),
  @(43282, 285838, 14299256, 571592151.58999908),
  @(43313, 286319, 16760893, 500307016.8599996),
  @(43344, 346191, 17475565, 390445461.99000031),
  @(43374, 368703, 18881092, 343255489.01999986),
  @(43405, 364995, 21480956, 283055203.18000036),
  @(43435, 413872, 23886894, 362189051.33999997),
  @(43466, 388698, 22251443, 527684797.20000058),
  @(43497, 361411, 19439654, 277577466.54000002),
  @(43525, 413709, 22641973, 380245193.7099995),
  @(43556, 375260, 23073264, 382720255.57000047),
  @(43586, 360580, 19812794, 403357922.25999963),
  @(43617, 340118, 15460735, 447721018.35999978),
  @(43647, 311274, 16310883, 575753661.58999991),
  @(43678, 314203, 19976419, 612311152.5999999),
  @(43709, 330376, 18386603, 367370258.88999987),
  @(43739, 372296, 21756458, 347900546.13999987),
  @(43770, 378443, 22088537, 278966807.40999997),
  @(43800, 417569, 25184632, 390474615.07000041),
  @(43831, 421091, 27628620, 382109311.6099999),
  @(43862, 324848, 17810145, 236730239.50999987),
  @(43891, 18594, 588614, 15387612.459999997),
  @(43922, $null, $null, 5772604.0499999989),
  @(43952, 1164, 28812, 5886047.6499999985),
  @(43983, 24096, 1012103, 41043876.279999994),
  @(44013, 77880, 4547189, 131664507.99000007),
  @(44044, 113132, 6900404, 187806238.88999993),
  @(44075, 191944, 5119702, 113692983.67999987),
  @(44105, 161265, 3569657, 56452161.50999999),
  @(44136, 1096, 38669, 2540600.39),
  @(44166, 371, 1190, 1541999.06),
  @(44197, 184, 4312, 578297.44000000006),
  @(44228, 2199, 172016, 3078391.83),
  @(44256, 564, 44837, 1801875.2599999998),
  @(44287, 2775, 79313, 1403873.0899999996),
  @(44317, 74894, 1745998, 24269951.48),
  @(44348, 161797, 4405524, 76461329.570000052),
  @(44378, 196960, 8324032, 160340515.23000002),
  @(44409, 208331, 9802316, 176647934.15999997),
  @(44440, 231102, 8619963, 149489929.46999991),
  @(44470, 267633, 12267458, 171206872.56999993),
  @(44501, 273566, 11553468, 155770776.09999987),
  @(44531, 286275, 12424144, 167331403.20999992),
  @(44562, 274588, 8596079, 94075318.640000015),
  @(44593, 211534, 9515434, 119127267.94),
  @(44621, 259459, 13139227, 161709307.42999998),
  @(44652, 279436, 17212360, 242821038),
  @(44682, 256859, 17689718, 295178224.75999993),
  @(44713, 238332, 19373645, 358897175.70999998),
  @(44743, 223495, 22480067, 437971136.07999998),
  @(44774, 207079, 22570646, 346544222.51999992),
  @(44805, 240886, 16968493, 263023175.62000003),
  @(44835, 282779, 20136835, 271289324.19000006),
  @(44866, 270510, 17451343, 210792609.53000003),
  @(44896, 296500, 19181267, 255170204.12999997)
)

$r = 2
foreach ($row in $data) {
  $wsNew.Cells.Item($r, 1).Value = $row[0]
  if ($row[1] -ne $null) { $wsNew.Cells.Item($r, 2).Value = $row[1] }
  if ($row[2] -ne $null) { $wsNew.Cells.Item($r, 3).Value = $row[2] }
  if ($row[3] -ne $null) { $wsNew.Cells.Item($r, 4).Value = $row[3] }
  $r = $r + 1
}

# --- Date formatting for column A (numFmtId 17 == "mmm-yy") ---
$wsNew.Range("A2:A61").NumberFormat = "mmm-yy"

# --- Final selection / active sheet state ---
$wsNew.Activate()
$wsNew.Range("F25").Select()
